$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 110-120 with refreshed sample data (values pulled in during
# this batch of the pipeline run).
$ws.Range("A110").Value = 58
$ws.Range("B110").Value = 2
$ws.Range("D110").Value = 1
$ws.Range("G110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("W110").Value = 1
$ws.Range("Y110").Value = 1
$ws.Range("Z110").Value = 0
$ws.Range("AU110").Value = 700
$ws.Range("A111").Value = 102.34
$ws.Range("B111").Value = 4
$ws.Range("H111").Value = 1
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 0
$ws.Range("R111").Value = 1
$ws.Range("S111").Value = 1
$ws.Range("W111").Value = 1
$ws.Range("AU111").Value = 1483.93
$ws.Range("A112").Value = 225.09
$ws.Range("B112").Value = 6
$ws.Range("D112").Value = 1
$ws.Range("I112").Value = 1
$ws.Range("N112").Value = 1
$ws.Range("Q112").Value = 1
$ws.Range("S112").Value = 1
$ws.Range("U112").Value = 1
$ws.Range("W112").Value = 1
$ws.Range("Z112").Value = 0
$ws.Range("AD112").Value = 1
$ws.Range("AU112").Value = 3601.44
$ws.Range("A113").Value = 26
$ws.Range("B113").Value = 1
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("R113").Value = 1
$ws.Range("AU113").Value = 415
$ws.Range("A114").Value = 47
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 1
$ws.Range("J114").Value = 0
$ws.Range("P114").Value = 1
$ws.Range("R114").Value = 1
$ws.Range("W114").Value = 1
$ws.Range("Z114").Value = 0
$ws.Range("AB114").Value = 1
$ws.Range("AU114").Value = 700
$ws.Range("A115").Value = 60
$ws.Range("B115").Value = 2
$ws.Range("C115").Value = 0
$ws.Range("O115").Value = 0
$ws.Range("R115").Value = 0
$ws.Range("W115").Value = 0
$ws.Range("AU115").Value = 1125
$ws.Range("A116").Value = 123.98
$ws.Range("B116").Value = 3
$ws.Range("E116").Value = 0
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 1
$ws.Range("L116").Value = 1
$ws.Range("M116").Value = 0
$ws.Range("N116").Value = 1
$ws.Range("P116").Value = 0
$ws.Range("S116").Value = 1
$ws.Range("U116").Value = 1
$ws.Range("V116").Value = 0
$ws.Range("AD116").Value = 1
$ws.Range("AH116").Value = 0
$ws.Range("AU116").Value = 1797.71
$ws.Range("A117").Value = 218.16
$ws.Range("B117").Value = 6
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("N117").Value = 1
$ws.Range("R117").Value = 0
$ws.Range("S117").Value = 1
$ws.Range("T117").Value = 0
$ws.Range("AD117").Value = 1
$ws.Range("AE117").Value = 0
$ws.Range("AU117").Value = 3163.32
$ws.Range("A118").Value = 37
$ws.Range("B118").Value = 1.5
$ws.Range("R118").Value = 1
$ws.Range("Y118").Value = 1
$ws.Range("Z118").Value = 0
$ws.Range("AU118").Value = 950
$ws.Range("A119").Value = 21
$ws.Range("B119").Value = 1
$ws.Range("E119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("M119").Value = 0
$ws.Range("P119").Value = 0
$ws.Range("R119").Value = 1
$ws.Range("U119").Value = 0
$ws.Range("W119").Value = 0
$ws.Range("Z119").Value = 1
$ws.Range("AM119").Value = 0
$ws.Range("AU119").Value = 649
$ws.Range("A120").Value = 178.29
$ws.Range("B120").Value = 5
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 1
$ws.Range("F120").Value = 1
$ws.Range("G120").Value = 0
$ws.Range("I120").Value = 1
$ws.Range("L120").Value = 1
$ws.Range("N120").Value = 1
$ws.Range("O120").Value = 0
$ws.Range("P120").Value = 1
$ws.Range("S120").Value = 1
$ws.Range("T120").Value = 1
$ws.Range("U120").Value = 1
$ws.Range("V120").Value = 0
$ws.Range("AD120").Value = 1
$ws.Range("AE120").Value = 0
$ws.Range("AU120").Value = 3030.93

# The refreshed batch has fewer rows than before; drop the now-stale
# trailing rows (previously 121-123) so the sheet's used range shrinks
# to A1:AU120 to match the new sample count.
$ws.Rows("121:123").Delete()
